# Add new branch pseudo-instructions (ble/bgt/blt/beqz/bnez) to the truth table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the existing "blez/bgtz/bltz" rows (28-30) to the new
#    "ble/bgt/blt" pseudo-ops - same opcode/funct, just a new mnemonic.
$ws.Range("B28").Value = "ble"
$ws.Range("B29").Value = "bgt"
$ws.Range("B30").Value = "blt"

# 2. Insert 5 fresh rows before row 31 (push the old "j"/"jal" rows from
#    31-32 down to 36-37) to make room for beqz, bnez, blez, bgtz, bltz.
$ws.Rows.Item(31).Resize(5).Insert()

# 3. beqz - row 31 (same shape as beq, different opcode/mnemonic)
$ws.Range("B31").Value = "beqz"
$ws.Range("C31").Value = "0x04"
$ws.Range("D31").Value = "X"
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "X"
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = "110_011"
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = "X"
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = "X"

# 4. bnez - row 32
$ws.Range("B32").Value = "bnez"
$ws.Range("C32").Value = "0x05"
$ws.Range("D32").Value = "X"
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = "X"
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = "110_001"
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = "X"
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = "X"

# 5. blez - row 33
$ws.Range("B33").Value = "blez"
$ws.Range("C33").Value = "0x06"
$ws.Range("D33").Value = "X"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = "X"
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = "111_101"
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = "X"
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = "X"

# 6. bgtz - row 34
$ws.Range("B34").Value = "bgtz"
$ws.Range("C34").Value = "0x07"
$ws.Range("D34").Value = "X"
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = "X"
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = "111_111"
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = "X"
$ws.Range("N34").Value = 1
$ws.Range("O34").Value = "X"

# 7. bltz - row 35
$ws.Range("B35").Value = "bltz"
$ws.Range("C35").Value = "0x01"
$ws.Range("D35").Value = "X"
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = "X"
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = "111_011"
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = "X"
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = "X"

# 8. Extend the "I" type merged label cell down over the new rows
#    (was A17:A30, now covers through the new rows, A17:A35).
$ws.Range("A17:A30").UnMerge()
$ws.Range("A17:A35").Merge()

# 9. Fix up the J-type opcode/funct placeholder cells that diff'd when the
#    two rows shifted from 31-32 to 36-37 (D/I on the "j" row, and the
#    style/string on a couple of the "jal" row cells).
$ws.Range("D36").Value = "X"
$ws.Range("I36").Value = "X"

# 10. Move the selection / scrolled view to roughly where the edit happened,
#     matching the saved view state after the edit.
$ws.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 13
